# Weekly update: insert a new price record at row 131, pushing the
# existing rows 131-263 down to 132-264. The new row starts as a copy
# of the (old) row 131 and only its date (column D) is changed to the
# new reporting date (2022-01-25 -> serial 44586).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 131..263 down by one, leaving a blank row 131.
$ws.Rows(131).Insert()

# Populate the new row 131 with a copy of what is now row 132 (the
# original row 131 data), then overwrite its date.
$ws.Range("A132:R132").Copy($ws.Range("A131"))
$ws.Range("D131").Value = 44586
